# Regenerate save_data to use K (strikeouts) instead of Strike# for column G,
# recalculating/writing the s_vals (strikeout values) for each outing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values computed from the regenerated save_data (replacing old Strike# column).
$kValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 3
    11 = 0
    12 = 2
    13 = 3
    14 = 4
    15 = 0
    16 = 4
    17 = 1
    18 = 1
    19 = 4
    20 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
